$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows 2-11 (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$data = @(
    @(4, 1, 5, 5),
    @(5, 1, 10, 10),
    @(8, 2, 5, 5),
    @(9, 2, 10, 10),
    @(1, 3, 5, 5),
    @(3, 3, 10, 10),
    @(6, 3, 15, 15),
    @(7, 3, 20, 20),
    @(2, 4, 5, 5),
    @(10, 4, 10, 10)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
